# refrigerators_update.xlsx — "hide notes in contents screen"
#
# Adds a new `hideInContents` column to the survey sheet (so the form
# renderer knows to hide "note" rows on the contents screen), then leaves
# the workbook's active sheet/selection the way the author's last save did.

$wb = $excel.ActiveWorkbook

# --- survey sheet: add the hideInContents column ----------------------
$survey = $wb.Worksheets.Item("survey")

# New header cell F1; this also appends the shared string and extends the
# sheet's used range/dimension to F4 automatically.
$survey.Range("F1").Value = "hideInContents"

# Column F is sized to fit its header text, same as the other survey
# columns (all of which are custom-width).
$survey.Columns.Item(6).ColumnWidth = 12.4987

# Cursor on the survey sheet ends up on the freshly added header cell.
$survey.Range("F2").Select() | Out-Null

# --- make "settings" the active/selected sheet -------------------------
$settings = $wb.Worksheets.Item("settings")
$settings.Activate() | Out-Null
